$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.399.76"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "1.627.82"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "0.9999"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "304.05"
$ws.Range("E6").Value = "  -1.45%  "
$ws.Range("D7").Value = "0.3792"
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "52.10"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "0.3631"
$ws.Range("D10").Value = "0.08106"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "1.226"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("D12").Value = "1.000"
$ws.Range("E12").Value = "  -0.19%  "
$ws.Range("D13").Value = "22.66"
$ws.Range("E13").Value = "  -1.83%  "
$ws.Range("D14").Value = "6.550"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("D16").Value = "7.224"
$ws.Range("E16").Value = "  -2.92%  "
$ws.Range("D17").Value = "1.619.76"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").Value = "93.61"
$ws.Range("E18").Value = "  -1.20%  "
$ws.Range("D19").Value = "0.06906"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").Value = "17.89"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "6.412"
$ws.Range("E22").Value = "  -2.39%  "
$ws.Range("D23").Value = "23.409.12"
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "3.221"
$ws.Range("E25").Value = "  +3.84%  "
$ws.Range("D26").Value = "2.429"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "21.14"
$ws.Range("E27").Value = "  -1.06%  "
$ws.Range("D28").Value = "149.23"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("D29").Value = "5.289"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").Value = "134.53"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "2.300"
$ws.Range("E31").Value = "  -5.17%  "
$ws.Range("D32").Value = "1.800.68"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "6.795"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").Value = "11.02"
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("D35").Value = "0.9553"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").Value = "0.02791"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "0.08835"
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "0.07201"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("D40").Value = "6.106"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").Value = "0.7083"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("D42").Value = "1.357"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").Value = "16.28"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "12.32"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").Value = "0.6480"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").Value = "2.328"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "0.9993"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("D49").Value = "0.07996"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "1.205"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").Value = "125.60"
$ws.Range("E51").Value = "  -4.10%  "
